$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-04-20"

# Update the April row label text (shared string)
$ws.Range("A5").Value = "April (through 04-20)"

# January 2022 (row 2, column I)
$ws.Range("I2").Value = 161

# April 2022 row (row 5) updates
$ws.Range("C5").Value = 21
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 33
$ws.Range("G5").Value = 43
$ws.Range("H5").Value = 74
$ws.Range("I5").Value = 87

# Total row (row 6) updates
$ws.Range("C6").Value = 149
$ws.Range("E6").Value = 232
$ws.Range("F6").Value = 143
$ws.Range("G6").Value = 241
$ws.Range("H6").Value = 497
$ws.Range("I6").Value = 523
